$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPath = "/sps/lsst/groups/auxtel/data/hack_usdf/my_postisrccd_img_forspectractor_2023/empty~holo4_003/20230119"

for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 4).Value = $newPath
}
